$ws = $excel.ActiveWorkbook.ActiveSheet
for ($r = 2; $r -le 41; $r++) {
    $cell = $ws.Cells.Item($r, 11)
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 11
    $cell.Font.Color = 0
}
